$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# --- Rename header row (row 1) from the "_old" / "_new" suffix convention to
# the explicit format-version suffixes "_FV2310" (old/base) and "_FV2404" (new/target).

$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("K1").Value = "diff"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- Turn the data range into a proper Excel table ("Table1") with an
# autofilter on the header row, spanning the whole used range A1:U81.

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U81"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false

# --- Freeze the header row (row 1) so it stays visible while scrolling.

$ws.Activate()
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
